$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
$ws.Range("D2").Value = "январь"
$ws.Range("E2").Value = "Месяц 2022 года"

# --- Well drilling date range updates ---
$ws.Range("C11").Value = "31.12.2021 19:00"
$ws.Range("D11").Value = "07.02.2022 19:00"

# --- Daily consumption values for row 11 (columns E:AI) ---
$ws.Range("E11").Value = 10509.839999999998
$ws.Range("F11").Value = 11485.44
$ws.Range("G11").Value = 10653.839999999998
$ws.Range("H11").Value = 23938.56
$ws.Range("I11").Value = 26178.479999999992
$ws.Range("J11").Value = 29482.559999999998
$ws.Range("K11").Value = 15465.599999999999
$ws.Range("L11").Value = 8470.800000000001
$ws.Range("M11").Value = 7762.32
$ws.Range("N11").Value = 9088.560000000001
$ws.Range("O11").Value = 17064.719999999998
$ws.Range("P11").Value = 16463.519999999997
$ws.Range("Q11").Value = 19965.6
$ws.Range("R11").Value = 21417.840000000004
$ws.Range("S11").Value = 17626.319999999996
$ws.Range("T11").Value = 9738
$ws.Range("U11").Value = 11371.68
$ws.Range("V11").Value = 10551.600000000002
$ws.Range("W11").Value = 8655.12
$ws.Range("X11").Value = 8474.399999999998
$ws.Range("Y11").Value = 7521.119999999999
$ws.Range("Z11").Value = 32140.080000000005
$ws.Range("AA11").Value = 10396.800000000001
$ws.Range("AB11").Value = 8688.240000000002
$ws.Range("AC11").Value = 24072.479999999992
$ws.Range("AD11").Value = 30810.24000000001
$ws.Range("AE11").Value = 23374.079999999998
$ws.Range("AF11").Value = 13427.280000000002
$ws.Range("AG11").Value = 10779.12
$ws.Range("AH11").Value = 9489.6
$ws.Range("AI11").Value = 2134.08
